# Add a new "age" variable row to the Variables dictionary sheet so that
# multiple timepoints can be represented (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a brand new blank row above the current row 3 (the "Horvath" row),
# shifting the existing data rows (Horvath..PedBE) down by one row.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the "age" variable definition.
$ws.Range("A3").Value = "age"
$ws.Range("B3").Value = "integer"
$ws.Range("C3").Value = "numeric"
$ws.Range("D3").Value = "Age when the clocks where measured"

# Match the formatting of the header-like rows above (row 2) for the new
# row's first four columns, same as Excel does when a user fills in a
# freshly inserted row to match its neighbours.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column D needs to grow since the new label is the longest text in the
# column now ("Age when the clocks where measured").
$ws.Columns.Item(4).ColumnWidth = 28.498697916666668

# The Variables sheet becomes the active/selected tab, with the selected
# cell parked at B14.
$ws.Activate()
$ws.Range("B14").Select()

Write-Host "Inserted age row and updated view state"
